$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 for a "Docentes responsáveis:" entry. This
# pushes "Programa resumido:" and everything below it down by one row.
# Excel's default Insert() copies formatting (including row height) down from
# the row above, which reproduces the row-height pattern seen in the target
# workbook (e.g. the "Short syllabus:" / "Syllabus:" label rows pick up the
# height of the row that used to sit above them).
$ws.Rows(12).Insert()

# New row 12: label only, no body text in B/C.
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13 is the old "Programa resumido:" row (shifted down); its body text
# changes from the professor name to "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 10 ("Objetivos:") body text is replaced with the professor reference.
$ws.Range("B10").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C10").Value = "2166002 - Sandra Giacomin Schneider"

# Row 15 ("Programa:", previously row 14) body text becomes "01/01/2012".
# Assigning that string straight to .Value would make Excel parse it as a
# date serial (and mint a new number-format style), so instead copy the
# already-existing literal text value from the "Ativação:" row (B8/C8, which
# already holds that exact string as text) using a values-only paste. That
# carries the text over as a plain string without disturbing the destination
# cell's existing style.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 18 ("Método:", previously row 17) body text is replaced.
$ws.Range("B18").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C18").Value = "2166002 - Sandra Giacomin Schneider"
